# Update countries & provincias Spain
# - Update the "datos actualizados" timestamp in A1
# - Refresh case counters for a handful of countries whose rows shifted
#   rank order as new data came in (country label + B:H numeric columns)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 19:35"

# --- Helper: write a whole data row (country name + 7 metric columns) -
function Set-CountryRow($row, $name, $total, $nuevos, $activos, $recuperados, $criticos, $muertesHoy, $muertes) {
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $total
    $ws.Cells.Item($row, 3).Value = $nuevos
    $ws.Cells.Item($row, 4).Value = $activos
    $ws.Cells.Item($row, 5).Value = $recuperados
    $ws.Cells.Item($row, 6).Value = $criticos
    $ws.Cells.Item($row, 7).Value = $muertesHoy
    $ws.Cells.Item($row, 8).Value = $muertes
}

# --- Rows whose totals were refreshed but keep their country label ----
Set-CountryRow 4   "Estados Unidos" 1735296 10021 481988 1151975 0 761 101333
Set-CountryRow 5   "Brasil"          396166  3806 158593  212827 0 197  24746
Set-CountryRow 12  "Turquia"         159797  1035 122793   32573 0  34   4431
Set-CountryRow 13  "India"           157935  7142  67452   85955 0 184   4528
Set-CountryRow 42  "Israel"           16793    36  14570    1942 0   0    281
Set-CountryRow 54  "Kazajistan"        9304   335   4768    4499 0   0     37
Set-CountryRow 94  "Somalia"           1731    20    265    1399 0   0     67
Set-CountryRow 168 "Siria"              121     0     43      74 0   0      4
Set-CountryRow 175 "Comoras"             87     0     24      61 0   1      2

# --- Rows that swapped rank with a neighbour: label AND totals change -
Set-CountryRow 150 "Suazilandia"            272  11 168 102 0 0  2
Set-CountryRow 151 "Mauritania"             268   0  15 240 0 0 13
Set-CountryRow 152 "Liberia"                266   0 144  95 0 1 27

Set-CountryRow 199 "Santa Lucia"             18   0  18   0 0 0  0
Set-CountryRow 201 "Belice"                  18   0  16   0 0 0  2

Set-CountryRow 205 "Groenlandia"             13   1  11   2 0 0  0
Set-CountryRow 206 "Islas Malvinas"          13   0  13   0 0 0  0
Set-CountryRow 207 "Santa Sede"              12   0   2  10 0 0  0
Set-CountryRow 208 "Islas Turcas y Caicos"   12   0  10   1 0 0  1
